$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24: new date, hours and activity text; formula already present and will auto-extend
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A24").Value = 43849
$ws.Range("B24").Value = 1.5
$ws.Range("D24").Value = "Wöchentliches Meeting + Vorbereitung"

# Row 25: new date, hours and activity text
$ws.Range("A23").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A25").Value = 43851
$ws.Range("B25").Value = 6
$ws.Range("D25").Value = "Einrichtung Teamviewer + Klassendesignmeeting"

# Move the selection like the author's cursor ended up at A26
$ws.Range("A26").Select()

$wb.Save()
